# Update market-price derived columns across multiple sheets
# (mirrors a scheduled scrape/recalculation of Teamcraft-style profit sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 955.875
$ws.Range("I18").Value = 955.875
$ws.Range("K18").Value = 955.875
$ws.Range("M18").Value = -671.875

$ws.Range("H29").Value = 300
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H38").Value = 696.0714
$ws.Range("I38").Value = 249.61539
$ws.Range("J38").Value = 6500
$ws.Range("K38").Value = 748.84617
$ws.Range("L38").Value = 19500
$ws.Range("M38").Value = -376.84617
$ws.Range("N38").Value = -20244

$ws.Range("H40").Value = 4988.6665
$ws.Range("J40").Value = 4988.6665
$ws.Range("L40").Value = 4988.6665
$ws.Range("N40").Value = -5338.6665

$ws.Range("H62").Value = 4500
$ws.Range("I62").Value = 4500
$ws.Range("K62").Value = 4500
$ws.Range("M62").Value = -3876

$ws.Range("H65").Value = 4500
$ws.Range("I65").Value = 4500
$ws.Range("K65").Value = 22500
$ws.Range("M65").Value = -19380

$ws.Range("H92").Value = 1038.7812
$ws.Range("I92").Value = 695.8261
$ws.Range("J92").Value = 1915.2222
$ws.Range("K92").Value = 695.8261
$ws.Range("L92").Value = 1915.2222
$ws.Range("M92").Value = 552.1739
$ws.Range("N92").Value = -4411.2222

$ws.Range("H104").Value = 546.7778
$ws.Range("I104").Value = 457.2857
$ws.Range("K104").Value = 1371.8571
$ws.Range("M104").Value = 375.1428999999998

$ws.Range("H127").Value = 1297.2632
$ws.Range("I127").Value = 837.4167
$ws.Range("J127").Value = 2085.5715
$ws.Range("K127").Value = 2512.2501
$ws.Range("L127").Value = 6256.7145
$ws.Range("M127").Value = 2447.7499
$ws.Range("N127").Value = -16176.7145

$ws.Range("H137").Value = 3471.3333
$ws.Range("I137").Value = 4654.4
$ws.Range("J137").Value = 1992.5
$ws.Range("K137").Value = 13963.2
$ws.Range("L137").Value = 5977.5
$ws.Range("M137").Value = -11413.2
$ws.Range("N137").Value = -11077.5

$ws.Range("H138").Value = 124965.58
$ws.Range("J138").Value = 144512.69
$ws.Range("L138").Value = 433538.07
$ws.Range("N138").Value = -443818.07

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21446.281
$ws.Range("I32").Value = 19767.703
$ws.Range("J32").Value = 52500
$ws.Range("K32").Value = 19767.703
$ws.Range("L32").Value = 52500
$ws.Range("M32").Value = -19480.703
$ws.Range("N32").Value = -53074

$ws.Range("H63").Value = 3354.4
$ws.Range("I63").Value = 3482.8333
$ws.Range("J63").Value = 2840.6667
$ws.Range("K63").Value = 3482.8333
$ws.Range("L63").Value = 2840.6667
$ws.Range("M63").Value = -2796.8333
$ws.Range("N63").Value = -4212.6667

$ws.Range("H66").Value = 3354.4
$ws.Range("I66").Value = 3482.8333
$ws.Range("J66").Value = 2840.6667
$ws.Range("K66").Value = 17414.1665
$ws.Range("L66").Value = 14203.3335
$ws.Range("M66").Value = -13982.1665
$ws.Range("N66").Value = -21067.3335

$ws.Range("H74").Value = 4084004.8
$ws.Range("I74").Value = 4763512
$ws.Range("K74").Value = 4763512
$ws.Range("M74").Value = -4762638

$ws.Range("H77").Value = 4084004.8
$ws.Range("I77").Value = 4763512
$ws.Range("K77").Value = 23817560
$ws.Range("M77").Value = -23813192

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2299.125
$ws.Range("I20").Value = 1723.5
$ws.Range("J20").Value = 2874.75
$ws.Range("K20").Value = 1723.5
$ws.Range("L20").Value = 2874.75
$ws.Range("M20").Value = -1476.5
$ws.Range("N20").Value = -3368.75

$ws.Range("H82").Value = 16250
$ws.Range("J82").Value = 60283
$ws.Range("L82").Value = 60283
$ws.Range("N82").Value = -61049

$ws.Range("H85").Value = 16250
$ws.Range("J85").Value = 60283
$ws.Range("L85").Value = 60283
$ws.Range("N85").Value = -62935

$ws.Range("H94").Value = 1844.2084
$ws.Range("I94").Value = 1837.4348
$ws.Range("K94").Value = 1837.4348
$ws.Range("M94").Value = -1386.4348

$ws.Range("H99").Value = 3469.3044
$ws.Range("I99").Value = 2670.647
$ws.Range("J99").Value = 5732.1665
$ws.Range("K99").Value = 2670.647
$ws.Range("L99").Value = 5732.1665
$ws.Range("M99").Value = -1172.647
$ws.Range("N99").Value = -8728.166499999999

$ws.Range("H105").Value = 3248.6377
$ws.Range("I105").Value = 3055.082
$ws.Range("J105").Value = 4724.5
$ws.Range("K105").Value = 3055.082
$ws.Range("L105").Value = 4724.5
$ws.Range("M105").Value = -1308.082
$ws.Range("N105").Value = -8218.5

$ws.Range("H134").Value = 3231.7778
$ws.Range("I134").Value = 2871.476
$ws.Range("K134").Value = 8614.428
$ws.Range("M134").Value = -6079.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 74999
$ws.Range("J81").Value = 74999
$ws.Range("L81").Value = 74999
$ws.Range("N81").Value = -76995

$ws.Range("H84").Value = 74999
$ws.Range("J84").Value = 74999
$ws.Range("L84").Value = 224997
$ws.Range("N84").Value = -234981

$ws.Range("H87").Value = 169499.5
$ws.Range("J87").Value = 169499.5
$ws.Range("L87").Value = 169499.5
$ws.Range("N87").Value = -171871.5

$ws.Range("H90").Value = 169499.5
$ws.Range("J90").Value = 169499.5
$ws.Range("L90").Value = 508498.5
$ws.Range("N90").Value = -520354.5

$ws.Range("H96").Value = 59497.75
$ws.Range("J96").Value = 59497.75
$ws.Range("L96").Value = 59497.75
$ws.Range("N96").Value = -64989.75

$ws.Range("H99").Value = 6110.8887
$ws.Range("I99").Value = 6083.1665
$ws.Range("K99").Value = 6083.1665
$ws.Range("M99").Value = -4585.1665

$ws.Range("H102").Value = 89797
$ws.Range("J102").Value = 89797
$ws.Range("L102").Value = 89797
$ws.Range("N102").Value = -94665

$ws.Range("H126").Value = 6110.8887
$ws.Range("I126").Value = 6083.1665
$ws.Range("K126").Value = 18249.4995
$ws.Range("M126").Value = -15779.4995

$ws.Range("H134").Value = 3873.4412
$ws.Range("I134").Value = 1677.75
$ws.Range("K134").Value = 5033.25
$ws.Range("M134").Value = -2498.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5480132
$ws.Range("I4").Value = 7902788
$ws.Range("K4").Value = 23708364
$ws.Range("M4").Value = -23708252

$ws.Range("H80").Value = 6553.2856
$ws.Range("J80").Value = 6582.6665
$ws.Range("L80").Value = 19747.9995
$ws.Range("N80").Value = -21619.9995

$ws.Range("H83").Value = 6553.2856
$ws.Range("J83").Value = 6582.6665
$ws.Range("L83").Value = 59243.9985
$ws.Range("N83").Value = -68603.9985

$ws.Range("H117").Value = 4392.3335
$ws.Range("J117").Value = 4394
$ws.Range("L117").Value = 13182
$ws.Range("N117").Value = -20066

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1375.7059
$ws.Range("I97").Value = 901.4583
$ws.Range("K97").Value = 901.4583
$ws.Range("M97").Value = -405.4583

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 417.25
$ws.Range("I9").Value = 293
$ws.Range("K9").Value = 293
$ws.Range("M9").Value = -69

$ws.Range("H22").Value = 3024.6667
$ws.Range("I22").Value = 3899
$ws.Range("J22").Value = 2849.8
$ws.Range("K22").Value = 3899
$ws.Range("L22").Value = 2849.8
$ws.Range("M22").Value = -3604
$ws.Range("N22").Value = -3439.8

$ws.Range("H27").Value = 3024.6667
$ws.Range("I27").Value = 3899
$ws.Range("J27").Value = 2849.8
$ws.Range("K27").Value = 3899
$ws.Range("L27").Value = 2849.8
$ws.Range("M27").Value = -3792
$ws.Range("N27").Value = -3063.8
